# Sample_Prep_Helper_Template.xlsx update
# Commit message: "Condensing of scripts and better logging"
#
# Summary of functional changes (derived from the OOXML diff):
#  - Augment_Specs sheet: the "naw_cell_code" column/header and its cell codes
#    (CELL14/CELL18/CELL21/CELL22/CELL35) are renamed to a condensed
#    "cell_code" scheme (CNCE14/CNPL18/CNHH21/CNDP22/CNPR35).
#  - Weighting_Segments sheet: the last header ("naw_cell_code") is renamed
#    to match ("cell_code").
#  - Valid_DMA_Codes sheet: the two identical "Name" headers are disambiguated
#    to "Name1" and "Name2".
#  - Various sheets also had their active-cell selection moved (harmless
#    view-state, reproduced here for fidelity).

$wb = $excel.ActiveWorkbook

$wsAugment    = $wb.Worksheets.Item("Augment_Specs")
$wsSpCode     = $wb.Worksheets.Item("Sp_Code_Freqs")
$wsValidDma   = $wb.Worksheets.Item("Valid_DMA_Codes")
$wsWeighting  = $wb.Worksheets.Item("Weighting_Segments")

# --- Augment_Specs ---------------------------------------------------------
# Header rename: naw_cell_code -> cell_code
$wsAugment.Range("A1").Value = "cell_code"

# Cell-code value rename (same row order as before)
$wsAugment.Range("A2").Value = "CNCE14"
$wsAugment.Range("A3").Value = "CNPL18"
$wsAugment.Range("A4").Value = "CNHH21"
$wsAugment.Range("A5").Value = "CNDP22"
$wsAugment.Range("A6").Value = "CNPR35"

# --- Valid_DMA_Codes ---------------------------------------------------------
# Disambiguate the two "Name" headers
$wsValidDma.Range("A1").Value = "Name1"
$wsValidDma.Range("B1").Value = "Name2"

# --- Weighting_Segments ---------------------------------------------------------
# Header rename: naw_cell_code -> cell_code
$wsWeighting.Range("H1").Value = "cell_code"

# --- View-state: active cell selections ------------------------------------
$wsAugment.Activate()
$wsAugment.Range("B10").Select()

$wsSpCode.Activate()
$wsSpCode.Range("B13").Select()

$wsValidDma.Activate()
$wsValidDma.Range("A2").Select()

$wsWeighting.Activate()
$wsWeighting.Range("H2").Select()
